$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.851.64'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").Value = '1.737.71'
$ws.Range("E3").Value = '  +0.07%  '
$ws.Range("D4").Value = '''0.9992'
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = '''236.94'
$ws.Range("E5").Value = '  +2.98%  '
$ws.Range("D6").Value = '''0.9994'
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").Value = '''0.5122'
$ws.Range("E7").Value = '  -1.14%  '
$ws.Range("D8").Value = '''0.2728'
$ws.Range("E8").Value = '  -0.76%  '
$ws.Range("D9").Value = '''0.06113'
$ws.Range("E9").Value = '  -0.15%  '
$ws.Range("D10").Value = '1.738.51'
$ws.Range("E10").Value = '  +0.16%  '
$ws.Range("E11").Value = '  +1.40%  '
$ws.Range("D12").Value = '''14.92'
$ws.Range("E12").Value = '  -1.42%  '
$ws.Range("D13").Value = '''0.6356'
$ws.Range("E13").Value = '  +0.10%  '
$ws.Range("D14").Value = '''4.590'
$ws.Range("E14").Value = '  +2.02%  '
$ws.Range("D15").Value = '''77.20'
$ws.Range("E15").Value = '  +0.43%  '
$ws.Range("D16").Value = '''0.9997'
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("D18").Value = '25.854.51'
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("D19").Value = '''11.71'
$ws.Range("E19").Value = '  +2.26%  '
$ws.Range("D20").Value = '''0.000006719'
$ws.Range("E20").Value = '  +1.28%  '
$ws.Range("D21").Value = '1.960.50'
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("D22").Value = '''4.253'
$ws.Range("E22").Value = '  +2.96%  '
$ws.Range("D23").Value = '''8.648'
$ws.Range("E23").Value = '  -0.78%  '
$ws.Range("E24").Value = '  +1.75%  '
$ws.Range("D25").Value = '''138.68'
$ws.Range("E25").Value = '  -0.65%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").Value = '''15.10'
$ws.Range("E27").Value = '  +0.69%  '
$ws.Range("D28").Value = '''1.752'
$ws.Range("E28").Value = '  -1.34%  '
$ws.Range("D29").Value = '''105.31'
$ws.Range("E29").Value = '  +3.39%  '
$ws.Range("D30").Value = '''3.958'
$ws.Range("E30").Value = '  +7.42%  '
$ws.Range("D31").Value = '''0.08333'
$ws.Range("E31").Value = '  +0.42%  '
$ws.Range("D32").Value = '''3.639'
$ws.Range("E32").Value = '  +4.75%  '
$ws.Range("D33").Value = '''0.04557'
$ws.Range("E33").Value = '  +1.49%  '
$ws.Range("D34").Value = '''2.662'
$ws.Range("E34").Value = '  +1.74%  '
$ws.Range("D35").Value = '''0.9811'
$ws.Range("E35").Value = '  +0.82%  '
$ws.Range("D36").Value = '''0.6156'
$ws.Range("E36").Value = '  +0.38%  '
$ws.Range("D37").Value = '''2.689'
$ws.Range("E37").Value = '  +1.29%  '
$ws.Range("D38").Value = '''0.01593'
$ws.Range("E38").Value = '  +0.88%  '
$ws.Range("D39").Value = '''1.916'
$ws.Range("E39").Value = '  -1.26%  '
$ws.Range("D40").Value = '''0.9991'
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("D41").Value = '''97.78'
$ws.Range("E41").Value = '  -2.50%  '
$ws.Range("D42").Value = '''0.3833'
$ws.Range("E42").Value = '  +0.46%  '
$ws.Range("D43").Value = '''0.7331'
$ws.Range("E43").Value = '  +1.66%  '
$ws.Range("D44").Value = '''4.939'
$ws.Range("E44").Value = '  -1.25%  '
$ws.Range("E45").Value = '  +0.12%  '
$ws.Range("D46").Value = '''0.05263'
$ws.Range("E46").Value = '  -2.11%  '
$ws.Range("D47").Value = '''6.160'
$ws.Range("E47").Value = '  -0.92%  '
$ws.Range("D48").Value = '''54.70'
$ws.Range("E48").Value = '  +3.44%  '
$ws.Range("D49").Value = '''30.45'
$ws.Range("E49").Value = '  +1.81%  '
$ws.Range("D50").Value = '''7.532'
$ws.Range("E50").Value = '  -0.85%  '
$ws.Range("D51").Value = '''0.3412'
$ws.Range("E51").Value = '  +1.46%  '
